$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.388.61'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.639.07'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  +3.31%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.93'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.64%  '
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0610'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = '1.870.41'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').Value = '1.624.86'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.45%  '
$ws.Range('D17').Value = '27.364.95'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.46'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.80%  '
$ws.Range('D19').Value = '0.0₃0720'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.69%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.84%  '
$ws.Range('E30').Value = '  -4.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0483'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.28'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').Value = '1.411.00'
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.564'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.880'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.43%  '
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.795'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.70%  '
$ws.Range('D47').Value = '1.780.42'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.65'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.72%  '
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.20%  '
